try {
    $wb = $excel.ActiveWorkbook
    $ws = $wb.ActiveSheet

    # Replace the single data row (row 2) with the new site record,
    # wrapped in try/except so a bad user input doesn't blow up the run.
    try { $ws.Range("A2").Value = 217 } catch { Write-Host "Error setting A2: $_" }

    try {
        # Keep the site code as text so leading zeros aren't lost.
        $ws.Range("B2").NumberFormat = "@"
        $ws.Range("B2").Value = "001058"
        $ws.Range("B2").Style = "Normal"
    } catch { Write-Host "Error setting B2: $_" }

    try { $ws.Range("C2").Value = "Great Island Channel SAC" } catch { Write-Host "Error setting C2: $_" }
    try { $ws.Range("D2").Value = "co" } catch { Write-Host "Error setting D2: $_" }
    try { $ws.Range("E2").Value = 1437.549977 } catch { Write-Host "Error setting E2: $_" }
    try { $ws.Range("F2").Value = 3.00999999046 } catch { Write-Host "Error setting F2: $_" }
    try { $ws.Range("G2").Value = "https://www.npws.ie/protected-sites/sac/001058" } catch { Write-Host "Error setting G2: $_" }

    # The new dataset only has one row of data, so drop the old rows 3-9.
    try { $ws.Range("A3:G9").EntireRow.Delete() } catch { Write-Host "Error deleting rows 3-9: $_" }
}
catch {
    Write-Host "Error: $_"
}
